$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on the numeric-looking Price/Volume columns so Excel
# COM does not auto-convert the assigned strings into Double/Date values
# (these columns hold display strings like "26.154.61" / "  -1.43%  ").
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '26.154.61'
$ws.Range("E2").Value = '  -1.43%  '
$ws.Range("D3").Value = '1.653.70'
$ws.Range("E3").Value = '  -1.79%  '
$ws.Range("D4").Value = '1.005'
$ws.Range("E4").Value = '  +0.34%  '
$ws.Range("D5").Value = '218.41'
$ws.Range("D6").Value = '0.5208'
$ws.Range("E6").Value = '  -2.19%  '
$ws.Range("E7").Value = '  +0.32%  '
$ws.Range("D8").Value = '0.2667'
$ws.Range("E8").Value = '  -0.41%  '
$ws.Range("D9").Value = '0.06325'
$ws.Range("E9").Value = '  -1.60%  '
$ws.Range("D10").Value = '21.10'
$ws.Range("E10").Value = '  -1.72%  '
$ws.Range("D11").Value = '0.07735'
$ws.Range("E11").Value = '  -0.63%  '
$ws.Range("D12").Value = '4.439'
$ws.Range("E12").Value = '  -1.38%  '
$ws.Range("D13").Value = '1.648.25'
$ws.Range("E13").Value = '  -2.08%  '
$ws.Range("D14").Value = '1.880.39'
$ws.Range("E14").Value = '  -1.82%  '
$ws.Range("D15").Value = '0.5467'
$ws.Range("E15").Value = '  -2.73%  '
$ws.Range("D16").Value = '0.0₅8228'
$ws.Range("E16").Value = '  -2.26%  '
$ws.Range("D17").Value = '64.85'
$ws.Range("E17").Value = '  -1.80%  '
$ws.Range("D18").Value = '26.194.62'
$ws.Range("E18").Value = '  -1.43%  '
$ws.Range("E19").Value = '  +0.37%  '
$ws.Range("E20").Value = '  -3.08%  '
$ws.Range("D21").Value = '192.61'
$ws.Range("E21").Value = '  -0.98%  '
$ws.Range("E22").Value = '  -2.36%  '
$ws.Range("D23").Value = '6.100'
$ws.Range("E23").Value = '  -4.65%  '
$ws.Range("D24").Value = '1.008'
$ws.Range("E24").Value = '  +0.50%  '
$ws.Range("D25").Value = '137.14'
$ws.Range("E25").Value = '  -4.66%  '
$ws.Range("D26").Value = '0.1238'
$ws.Range("E26").Value = '  -2.56%  '
$ws.Range("E27").Value = '  -3.21%  '
$ws.Range("E28").Value = '  -0.66%  '
$ws.Range("D29").Value = '1.413'
$ws.Range("E29").Value = '  -0.15%  '
$ws.Range("D30").Value = '0.06036'
$ws.Range("E30").Value = '  -1.31%  '
$ws.Range("D31").Value = '1.282'
$ws.Range("E31").Value = '  +0.27%  '
$ws.Range("D32").Value = '3.569'
$ws.Range("D33").Value = '3.335'
$ws.Range("E33").Value = '  -3.71%  '
$ws.Range("D34").Value = '1.651'
$ws.Range("E34").Value = '  -2.82%  '
$ws.Range("D35").Value = '0.9809'
$ws.Range("E35").Value = '  -3.64%  '
$ws.Range("D36").Value = '2.411'
$ws.Range("E36").Value = '  -0.39%  '
$ws.Range("E37").Value = '  -0.88%  '
$ws.Range("D38").Value = '0.5947'
$ws.Range("E38").Value = '  +4.12%  '
$ws.Range("E39").Value = '  -2.94%  '
$ws.Range("D40").Value = '5.957'
$ws.Range("E40").Value = '  -0.33%  '
$ws.Range("D41").Value = '0.8634'
$ws.Range("E41").Value = '  -0.44%  '
$ws.Range("E42").Value = '  +0.20%  '
$ws.Range("D43").Value = '1.038.57'
$ws.Range("E43").Value = '  -1.65%  '
$ws.Range("D44").Value = '99.66'
$ws.Range("E44").Value = '  -0.50%  '
$ws.Range("D45").Value = '1.794.06'
$ws.Range("E45").Value = '  -2.26%  '
$ws.Range("B46").Value = 'BabyDogeCoin'
$ws.Range("C46").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D46").Value = '0.0₈112'
$ws.Range("E46").Value = '  -0.45%  '
$ws.Range("B47").Value = 'Aave'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D47").Value = '57.11'
$ws.Range("E47").Value = '  -0.13%  '
$ws.Range("E48").Value = '  +0.38%  '
$ws.Range("D49").Value = '8.109'
$ws.Range("E49").Value = '  -0.62%  '
$ws.Range("D50").Value = '0.05180'
$ws.Range("E50").Value = '  -0.56%  '
$ws.Range("E51").Value = '  +3.82%  '
